# Fix Training Data Issue (#48)
#
# The "Date" column (BF) on Sheet1 holds the game date for every row of
# team stats, but it was stored in the wrong/inconsistent textual form
# "5-16-2012-13" (an artifact of how NBA stats were shown - mixing the
# game day with the season label). Correct it to the real ISO-style game
# date "2013-05-16" for every data row (rows 2-31).
#
# The value must stay plain text (not be re-interpreted by Excel as a
# date serial number), and must end up with no visible/semantic style
# change on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$correctDate = "2013-05-16"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Range("BF$row")
    # A leading apostrophe forces Excel to store the assignment as literal
    # text instead of re-parsing "2013-05-16" into a date serial value.
    $cell.Value = "'" + $correctDate
}

# Restore the default "Normal" style on the whole column range so the
# text/quote-prefix formatting Excel applies while typing doesn't leave
# a visible style change behind on these cells.
$rangeAddress = "BF" + $firstRow + ":BF" + $lastRow
$ws.Range($rangeAddress).Style = "Normal"
